$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text columns (Coin name / Link) - safe as plain text ---
$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("B21").Value = "Uniswap"
$ws.Range("B25").Value = "Cosmos"
$ws.Range("B26").Value = "Stellar"
$ws.Range("B38").Value = "VeChain"
$ws.Range("B39").Value = "MXToken"
$ws.Range("B40").Value = "FraxShare"
$ws.Range("B41").Value = "Aave"
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("B43").Value = "RenderToken"
$ws.Range("B44").Value = "Quant"
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("B48").Value = "Aptos"
$ws.Range("B49").Value = "Decentraland"
$ws.Range("B51").Value = "Maker"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("C26").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("C49").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"

# --- Price column (D) - force text format so numeric-looking strings
#     ("1.001", "254.28", ...) are NOT coerced into Excel numbers, then
#     restore the default "Normal" style so no stray number-format is left
#     applied to the cells (matches the original un-styled cells). ---
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"
$ws.Range("D2").Value = "31.553.83"
$ws.Range("D3").Value = "1.992.50"
$ws.Range("D4").Value = "1.001"
$ws.Range("D5").Value = "0.8086"
$ws.Range("D6").Value = "254.28"
$ws.Range("D7").Value = "0.9993"
$ws.Range("D8").Value = "0.3435"
$ws.Range("D9").Value = "25.86"
$ws.Range("D10").Value = "0.07099"
$ws.Range("D11").Value = "0.8465"
$ws.Range("D12").Value = "0.08118"
$ws.Range("D13").Value = "102.52"
$ws.Range("D14").Value = "1.987.06"
$ws.Range("D15").Value = "5.496"
$ws.Range("D16").Value = "276.23"
$ws.Range("D17").Value = "31.538.84"
$ws.Range("D18").Value = "14.01"
$ws.Range("D19").Value = "0.000008009"
$ws.Range("D20").Value = "2.248.46"
$ws.Range("D21").Value = "5.709"
$ws.Range("D22").Value = "0.9990"
$ws.Range("D23").Value = "1.000"
$ws.Range("D24").Value = "6.922"
$ws.Range("D25").Value = "9.718"
$ws.Range("D26").Value = "0.1555"
$ws.Range("D27").Value = "165.92"
$ws.Range("D28").Value = "19.83"
$ws.Range("D29").Value = "2.221"
$ws.Range("D30").Value = "1.569"
$ws.Range("D31").Value = "1.361"
$ws.Range("D32").Value = "4.590"
$ws.Range("D33").Value = "4.334"
$ws.Range("D34").Value = "0.05197"
$ws.Range("D35").Value = "1.220"
$ws.Range("D36").Value = "0.7539"
$ws.Range("D37").Value = "2.808"
$ws.Range("D38").Value = "0.02011"
$ws.Range("D39").Value = "2.943"
$ws.Range("D40").Value = "6.665"
$ws.Range("D41").Value = "78.69"
$ws.Range("D42").Value = "0.4691"
$ws.Range("D43").Value = "2.084"
$ws.Range("D44").Value = "106.69"
$ws.Range("D45").Value = "0.8619"
$ws.Range("D46").Value = "0.9992"
$ws.Range("D47").Value = "9.978"
$ws.Range("D48").Value = "7.540"
$ws.Range("D49").Value = "0.4293"
$ws.Range("D50").Value = "36.59"
$ws.Range("D51").Value = "929.62"
$dRange.Style = "Normal"

# --- Volume(1h) column (E) - percentage text with padding spaces; ---
# --- already safe as text because of the surrounding whitespace,  ---
# --- but force text format too for robustness / symmetry.          ---
$eRange = $ws.Range("E2:E51")
$eRange.NumberFormat = "@"
$ws.Range("E2").Value = "  +4.06%  "
$ws.Range("E3").Value = "  +6.40%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("E5").Value = "  +72.31%  "
$ws.Range("E6").Value = "  +4.31%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  +19.49%  "
$ws.Range("E9").Value = "  +17.40%  "
$ws.Range("E10").Value = "  +10.18%  "
$ws.Range("E11").Value = "  +17.44%  "
$ws.Range("E12").Value = "  +4.52%  "
$ws.Range("E13").Value = "  +7.20%  "
$ws.Range("E14").Value = "  +5.94%  "
$ws.Range("E15").Value = "  +7.29%  "
$ws.Range("E16").Value = "  -1.13%  "
$ws.Range("E17").Value = "  +4.05%  "
$ws.Range("E18").Value = "  +7.91%  "
$ws.Range("E19").Value = "  +7.88%  "
$ws.Range("E20").Value = "  +5.85%  "
$ws.Range("E21").Value = "  +9.44%  "
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("E24").Value = "  +11.18%  "
$ws.Range("E25").Value = "  +7.46%  "
$ws.Range("E26").Value = "  +62.57%  "
$ws.Range("E27").Value = "  +1.65%  "
$ws.Range("E28").Value = "  +6.29%  "
$ws.Range("E29").Value = "  +18.29%  "
$ws.Range("E30").Value = "  +6.84%  "
$ws.Range("E31").Value = "  +3.05%  "
$ws.Range("E32").Value = "  +8.86%  "
$ws.Range("E33").Value = "  +5.76%  "
$ws.Range("E34").Value = "  +8.19%  "
$ws.Range("E36").Value = "  +9.81%  "
$ws.Range("E37").Value = "  +3.73%  "
$ws.Range("E38").Value = "  +7.63%  "
$ws.Range("E39").Value = "  +4.74%  "
$ws.Range("E40").Value = "  +7.17%  "
$ws.Range("E41").Value = "  +6.11%  "
$ws.Range("E42").Value = "  +10.86%  "
$ws.Range("E43").Value = "  +8.08%  "
$ws.Range("E44").Value = "  +6.01%  "
$ws.Range("E45").Value = "  +4.43%  "
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("E47").Value = "  +4.39%  "
$ws.Range("E48").Value = "  +9.05%  "
$ws.Range("E49").Value = "  +10.00%  "
$ws.Range("E50").Value = "  +4.32%  "
$ws.Range("E51").Value = "  +3.67%  "
$eRange.Style = "Normal"
